$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the taxon-observation data between row 14 and row 15
# (location/date/observer metadata columns are identical between the two
# rows and are left untouched).

# --- Row 14: adopt former row-15 values (Kolflarnlav / Carbonicola anthracophila) ---
$ws.Range("A14").Value = 111380345
$ws.Range("B14").Value = 77267
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6446
$ws.Range("F14").Value = "Kolflarnlav"
$ws.Range("G14").Value = "Carbonicola anthracophila"
$ws.Range("H14").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("I14").ClearContents()
$ws.Range("J14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("Q14").Value = 364908.1256513004
$ws.Range("R14").Value = 6872135.474104149
$ws.Range("AC14").Value = "Växer på gammal kolad tallstubbe i kontinuitetsskog"
$ws.Range("AI14").Value = "Tallskog. Kontinuitetsskog"
$ws.Range("AJ14").Value = "tall"
$ws.Range("AK14").Value = "Pinus sylvestris"
$ws.Range("AO14").Value = "Pinus sylvestris"

# --- Row 15: adopt former row-14 values (Plattlummer / Lycopodium complanatum) ---
$ws.Range("A15").Value = 111379229
$ws.Range("B15").Value = 95538
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 221941
$ws.Range("F15").Value = "Plattlummer"
$ws.Range("G15").Value = "Lycopodium complanatum"
$ws.Range("H15").Value = "L."
$ws.Range("I15").Value = "'10"
$ws.Range("J15").Value = "m²"
$ws.Range("L15").ClearContents()
$ws.Range("Q15").Value = 364945.755472637
$ws.Range("R15").Value = 6872251.713583581
$ws.Range("AC15").Value = "Plattlummer växer i k-skog"
$ws.Range("AI15").Value = "Barrblandskog med gamla tallar och senvuxna granar. Kontinuitetsskog"
$ws.Range("AJ15").ClearContents()
$ws.Range("AK15").ClearContents()
$ws.Range("AO15").ClearContents()
